{"js": "// Update the two-digit-division answer table: each <w:t> run in the\n// worked-example cells gets its quotient/remainder expression swapped for a\n// new one, per the commit's regenerated answer key. Matched in document\n// order so the handful of values that are reused as a later replacement's\n// target (e.g. \"35\u00f77=5, 0\") cannot be re-matched after they've already been\n// written.\nconst replacements = [\n  [\"36\u00f79=4, 0\", \"30\u00f76=5, 0\"],\n  [\"30\u00f72=15, 0\", \"23\u00f77=3, 2\"],\n  [\"72\u00f73=24, 0\", \"92\u00f72=46, 0\"],\n  [\"23\u00f79=2, 5\", \"26\u00f74=6, 2\"],\n  [\"24\u00f72=12, 0\", \"67\u00f74=16, 3\"],\n  [\"88\u00f75=17, 3\", \"52\u00f72=26, 0\"],\n  [\"33\u00f73=11, 0\", \"86\u00f75=17, 1\"],\n  [\"70\u00f79=7, 7\", \"53\u00f77=7, 4\"],\n  [\"84\u00f75=16, 4\", \"58\u00f78=7, 2\"],\n  [\"32\u00f79=3, 5\", \"61\u00f78=7, 5\"],\n  [\"66\u00f77=9, 3\", \"70\u00f77=10, 0\"],\n  [\"39\u00f74=9, 3\", \"37\u00f79=4, 1\"],\n  [\"87\u00f76=14, 3\", \"35\u00f72=17, 1\"],\n  [\"90\u00f75=18, 0\", \"83\u00f78=10, 3\"],\n  [\"53\u00f78=6, 5\", \"11\u00f77=1, 4\"],\n  [\"85\u00f74=21, 1\", \"62\u00f72=31, 0\"],\n  [\"93\u00f77=13, 2\", \"68\u00f78=8, 4\"],\n  [\"56\u00f78=7, 0\", \"52\u00f72=26, 0\"],\n  [\"71\u00f74=17, 3\", \"91\u00f73=30, 1\"],\n  [\"35\u00f77=5, 0\", \"88\u00f77=12, 4\"],\n  [\"44\u00f75=8, 4\", \"47\u00f75=9, 2\"],\n  [\"80\u00f77=11, 3\", \"84\u00f73=28, 0\"],\n  [\"29\u00f73=9, 2\", \"58\u00f75=11, 3\"],\n  [\"94\u00f79=10, 4\", \"68\u00f77=9, 5\"],\n  [\"79\u00f73=26, 1\", \"35\u00f77=5, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${before}\"`);\n  }\n\n  // Each \"before\" string is unique in the source document, so the first hit\n  // is the only (and correct) one.\n  results.items[0].insertText(after, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the two-digit-division answer table: each worked-example cell's\n# quotient/remainder expression is swapped for a new one, per the commit's\n# regenerated answer key. Pairs are listed (and applied) in document order so\n# that a value re-used later as a replacement target (e.g. \"35\u00f77=5, 0\", which\n# is the 20th cell's ORIGINAL text but becomes the 25th cell's NEW text)\n# can't be matched again after it has already been written.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"36\u00f79=4, 0\", \"30\u00f76=5, 0\"),\n  @(\"30\u00f72=15, 0\", \"23\u00f77=3, 2\"),\n  @(\"72\u00f73=24, 0\", \"92\u00f72=46, 0\"),\n  @(\"23\u00f79=2, 5\", \"26\u00f74=6, 2\"),\n  @(\"24\u00f72=12, 0\", \"67\u00f74=16, 3\"),\n  @(\"88\u00f75=17, 3\", \"52\u00f72=26, 0\"),\n  @(\"33\u00f73=11, 0\", \"86\u00f75=17, 1\"),\n  @(\"70\u00f79=7, 7\", \"53\u00f77=7, 4\"),\n  @(\"84\u00f75=16, 4\", \"58\u00f78=7, 2\"),\n  @(\"32\u00f79=3, 5\", \"61\u00f78=7, 5\"),\n  @(\"66\u00f77=9, 3\", \"70\u00f77=10, 0\"),\n  @(\"39\u00f74=9, 3\", \"37\u00f79=4, 1\"),\n  @(\"87\u00f76=14, 3\", \"35\u00f72=17, 1\"),\n  @(\"90\u00f75=18, 0\", \"83\u00f78=10, 3\"),\n  @(\"53\u00f78=6, 5\", \"11\u00f77=1, 4\"),\n  @(\"85\u00f74=21, 1\", \"62\u00f72=31, 0\"),\n  @(\"93\u00f77=13, 2\", \"68\u00f78=8, 4\"),\n  @(\"56\u00f78=7, 0\", \"52\u00f72=26, 0\"),\n  @(\"71\u00f74=17, 3\", \"91\u00f73=30, 1\"),\n  @(\"35\u00f77=5, 0\", \"88\u00f77=12, 4\"),\n  @(\"44\u00f75=8, 4\", \"47\u00f75=9, 2\"),\n  @(\"80\u00f77=11, 3\", \"84\u00f73=28, 0\"),\n  @(\"29\u00f73=9, 2\", \"58\u00f75=11, 3\"),\n  @(\"94\u00f79=10, 4\", \"68\u00f77=9, 5\"),\n  @(\"79\u00f73=26, 1\", \"35\u00f77=5, 0\")\n)\n\nforeach ($pair in $pairs) {\n  $findText = $pair[0]\n  $replaceText = $pair[1]\n\n  # Re-grab the whole-document range each time and search from the top;\n  # each \"before\" string is unique, and Replace=1 (wdReplaceOne) stops after\n  # the first (only) hit.\n  $r = $d.Content\n  $found = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n  if (-not $found) {\n    throw \"edit.ps1: no match found for '$findText'\"\n  }\n}\n"}
